# Auto-generated edit script applying the Excalibur_Profits diff
# Updates numeric values in columns H-N across rows on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2670.6
$ws.Range("I15").Value = 2670.6
$ws.Range("K15").Value = 8011.799999999999
$ws.Range("M15").Value = -7842.799999999999
$ws.Range("H74").Value = 6265.0713
$ws.Range("I74").Value = 4447.1665
$ws.Range("K74").Value = 4447.1665
$ws.Range("M74").Value = -3511.1665
$ws.Range("H77").Value = 6265.0713
$ws.Range("I77").Value = 4447.1665
$ws.Range("K77").Value = 22235.8325
$ws.Range("M77").Value = -17555.8325
$ws.Range("H86").Value = 1921.2858
$ws.Range("I86").Value = 1014
$ws.Range("K86").Value = 1014
$ws.Range("M86").Value = 109
$ws.Range("H89").Value = 1921.2858
$ws.Range("I89").Value = 1014
$ws.Range("K89").Value = 5070
$ws.Range("M89").Value = 546
$ws.Range("H100").Value = 5405.2856
$ws.Range("I100").Value = 2221.111
$ws.Range("K100").Value = 2221.111
$ws.Range("M100").Value = -1680.111
$ws.Range("H106").Value = 1820.625
$ws.Range("I106").Value = 1353.2
$ws.Range("J106").Value = 2599.6667
$ws.Range("K106").Value = 1353.2
$ws.Range("L106").Value = 2599.6667
$ws.Range("M106").Value = -722.2
$ws.Range("N106").Value = -3861.6667
$ws.Range("H112").Value = 6252.6226
$ws.Range("J112").Value = 6341.4424
$ws.Range("L112").Value = 19024.3272
$ws.Range("N112").Value = -21240.3272
$ws.Range("H132").Value = 77643.64999999999
$ws.Range("I132").Value = 84998.39
$ws.Range("K132").Value = 254995.17
$ws.Range("M132").Value = -252465.17
$ws.Range("H137").Value = 638960.5600000001
$ws.Range("I137").Value = 2270.1
$ws.Range("J137").Value = 1049728.6
$ws.Range("K137").Value = 6810.299999999999
$ws.Range("L137").Value = 3149185.8
$ws.Range("M137").Value = -4260.299999999999
$ws.Range("N137").Value = -3154285.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18055.672
$ws.Range("I32").Value = 19242
$ws.Range("J32").Value = 260.75
$ws.Range("K32").Value = 19242
$ws.Range("L32").Value = 260.75
$ws.Range("M32").Value = -18955
$ws.Range("N32").Value = -834.75
$ws.Range("H74").Value = 3064.4443
$ws.Range("I74").Value = 1082.8572
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 1082.8572
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -208.8571999999999
$ws.Range("N74").Value = -11748
$ws.Range("H77").Value = 3064.4443
$ws.Range("I77").Value = 1082.8572
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 5414.286
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -1046.286
$ws.Range("N77").Value = -58736
$ws.Range("H97").Value = 946.7143
$ws.Range("I97").Value = 513.913
$ws.Range("K97").Value = 513.913
$ws.Range("M97").Value = -17.91300000000001
$ws.Range("H132").Value = 716316.0600000001
$ws.Range("I132").Value = 771282.7
$ws.Range("K132").Value = 2313848.1
$ws.Range("M132").Value = -2311318.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1723.7858
$ws.Range("I20").Value = 1593.3
$ws.Range("K20").Value = 1593.3
$ws.Range("M20").Value = -1346.3
$ws.Range("H22").Value = 399.85715
$ws.Range("I22").Value = 399.85715
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 399.85715
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -226.85715
$ws.Range("N22").ClearContents()
$ws.Range("H105").Value = 3588.625
$ws.Range("I105").Value = 3588.625
$ws.Range("K105").Value = 3588.625
$ws.Range("M105").Value = -1841.625
$ws.Range("H134").Value = 547550.1
$ws.Range("I134").Value = 508045.72
$ws.Range("K134").Value = 1524137.16
$ws.Range("M134").Value = -1521602.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17232.1
$ws.Range("I31").Value = 6963
$ws.Range("J31").Value = 30660.924
$ws.Range("K31").Value = 6963
$ws.Range("L31").Value = 30660.924
$ws.Range("M31").Value = -6668
$ws.Range("N31").Value = -31250.924
$ws.Range("H34").Value = 17232.1
$ws.Range("I34").Value = 6963
$ws.Range("J34").Value = 30660.924
$ws.Range("K34").Value = 6963
$ws.Range("L34").Value = 30660.924
$ws.Range("M34").Value = -6761
$ws.Range("N34").Value = -31064.924
$ws.Range("H58").Value = 1773324.1
$ws.Range("I58").Value = 4121123
$ws.Range("J58").Value = 12475
$ws.Range("K58").Value = 4121123
$ws.Range("L58").Value = 12475
$ws.Range("M58").Value = -4120920
$ws.Range("N58").Value = -12881
$ws.Range("H134").Value = 8047.5454
$ws.Range("I134").Value = 8926.474
$ws.Range("J134").Value = 2481
$ws.Range("K134").Value = 26779.422
$ws.Range("L134").Value = 7443
$ws.Range("M134").Value = -24244.422
$ws.Range("N134").Value = -12513
$ws.Range("H136").Value = 1773324.1
$ws.Range("I136").Value = 4121123
$ws.Range("J136").Value = 12475
$ws.Range("K136").Value = 12363369
$ws.Range("L136").Value = 37425
$ws.Range("M136").Value = -12360819
$ws.Range("N136").Value = -42525

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 12641.556
$ws.Range("J110").Value = 25030
$ws.Range("L110").Value = 75090
$ws.Range("N110").Value = -83270
$ws.Range("H132").Value = 647.5
$ws.Range("I132").Value = 640
$ws.Range("K132").Value = 5760
$ws.Range("M132").Value = -3230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4984.154
$ws.Range("I70").Value = 4903.1665
$ws.Range("J70").Value = 5053.5713
$ws.Range("K70").Value = 4903.1665
$ws.Range("L70").Value = 5053.5713
$ws.Range("M70").Value = -4633.1665
$ws.Range("N70").Value = -5593.5713
$ws.Range("H73").Value = 4984.154
$ws.Range("I73").Value = 4903.1665
$ws.Range("J73").Value = 5053.5713
$ws.Range("K73").Value = 4903.1665
$ws.Range("L73").Value = 5053.5713
$ws.Range("M73").Value = -3967.1665
$ws.Range("N73").Value = -6925.5713
$ws.Range("H132").Value = 40486650
$ws.Range("I132").Value = 50606556
$ws.Range("J132").Value = 7024.6
$ws.Range("K132").Value = 151819668
$ws.Range("L132").Value = 21073.8
$ws.Range("M132").Value = -151817138
$ws.Range("N132").Value = -26133.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2399.6
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 1000
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 2399.6
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 5000
$ws.Range("N71").Value = -12488
$ws.Range("H82").Value = 1201
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 1502.5
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 1502.5
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -2224.5
$ws.Range("H85").Value = 1201
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 1502.5
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 1502.5
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -3998.5
$ws.Range("H132").Value = 1454703.2
$ws.Range("J132").Value = 5007.4546
$ws.Range("L132").Value = 15022.3638
$ws.Range("N132").Value = -20082.3638
$ws.Range("H136").Value = 3810.5625
$ws.Range("I136").Value = 2906.5454
$ws.Range("J136").Value = 5799.4
$ws.Range("K136").Value = 8719.636200000001
$ws.Range("L136").Value = 17398.2
$ws.Range("M136").Value = -6169.636200000001
$ws.Range("N136").Value = -22498.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2744.5
$ws.Range("I107").Value = 1961.5
$ws.Range("J107").Value = 3266.5
$ws.Range("K107").Value = 5884.5
$ws.Range("L107").Value = 9799.5
$ws.Range("M107").Value = -3964.5
$ws.Range("N107").Value = -13639.5
$ws.Range("H122").Value = 3374.6365
$ws.Range("I122").Value = 2990.875
$ws.Range("K122").Value = 8972.625
$ws.Range("M122").Value = -6522.625

Write-Output "Applied 202 cell updates across 8 sheets"